# Update cryptos list - GitHub Actions data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: the Price column (D) stores plain-text values (e.g. "30.389.86"
# using "." as a thousands separator, or fixed-precision decimals like
# "1.0000"); these must stay text, not be re-interpreted as numbers. Each
# written Price cell is briefly forced to Text format so the literal
# string is preserved exactly, then the format is cleared again so the
# cell's style matches the source file (plain, unstyled text).

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.389.86"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.73%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.63"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.92%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.34"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.77%  "

# Row 6 - USDC
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.03%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4772"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.87%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2876"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.12%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06516"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.33%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  -0.43%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07759"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.16%  "

# Row 12 - was WrappedEther, now Litecoin
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.60"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.62%  "

# Row 13 - was Litecoin, now WrappedEther
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.880.25"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.59%  "

# Row 14 - Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7355"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.94%  "

# Row 15 - Polkadot
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.127"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.97%  "

# Row 16 - BitcoinCash
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.33"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.53%  "

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.370.40"
$ws.Range("D17").ClearFormats()

# Row 18 - Avalanche
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.38"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.54%  "

# Row 19 - ShibaInu
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007538"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.32%  "

# Row 20 - Dai
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.11%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.127.46"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.21%  "

# Row 22 - BinanceUSD
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.06%  "

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.232"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.13%  "

# Row 24 - Chainlink
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.164"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.17%  "

# Row 25 - Cosmos
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.239"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.05%  "

# Row 26 - Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.51"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.46%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.90"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.50%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.957"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.80%  "

# Row 29 - Stellar
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09981"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.85%  "

# Row 30 - Toncoin: no change

# Row 31 - PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.511"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.64%  "

# Row 32 - Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.310"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.54%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.071"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.99%  "

# Row 34 - Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04748"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.74%  "

# Row 35 - ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.120"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.26%  "

# Row 36 - ImmutableX
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6959"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.01%  "

# Row 37 - HuobiToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.09%  "

# Row 38 - VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01858"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.08%  "

# Row 39 - MXToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.754"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.25%  "

# Row 40 - FraxShare
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.268"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.38%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8425"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.33%  "

# Row 42 - was Aave, now RenderToken
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.905"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.78%  "

# Row 43 - was RenderToken, now TheSandbox
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4162"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.99%  "

# Row 44 - was TheSandbox, now Aave
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.24"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.70%  "

# Row 45 - PaxDollar
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9995"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.06%  "

# Row 46 - Quant
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.86"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.72%  "

# Row 47 - Aptos
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.087"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.05%  "

# Row 48 - EnergySwap
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.200"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.35%  "

# Row 49 - Elrond
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.11"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.62%  "

# Row 50 - Maker
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "912.95"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.67%  "

# Row 51 - Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05592"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.70%  "
